$p = $ppt.ActivePresentation

# 1. Delete the "Key Feature: Topics" slide.
#    Locate it by its title text rather than assuming a fixed index, then
#    remove it from the deck (PowerPoint shifts everything after it up).
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    $isTarget = $false
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "Key Feature: Topics") {
                $isTarget = $true
            }
        }
    }
    if ($isTarget) {
        $s.Delete()
    }
}

# 2. Update the cached date field text (datetimeFigureOut) from 7/16/2012 to
#    7/22/2012 across the slide master and all slide layouts.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "7/16/2012") {
        $sh.TextFrame.TextRange.Text = "7/22/2012"
    }
}

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "7/16/2012") {
            $sh.TextFrame.TextRange.Text = "7/22/2012"
        }
    }
}
